$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1; this shifts the existing header (row 1) and
# all data rows down by one.
$ws.Rows.Item(1).Insert()

# New row 1: a note describing the accepted date-header formats.
$ws.Range("A1").Value = "Note: The date header (Row 2) supports: '2023 Annual', '2023 Q1', '2023-01'"
$ws.Range("A1").Font.Italic = $true
$ws.Range("A1").Font.Color = 255

# Merge A1:D1 so the note spans the full table width.
$ws.Range("A1:D1").Merge()

# Row 2 (previously row 1) keeps the bold/centered header style; update the
# year labels to include the "Annual" period qualifier.
$ws.Range("B2").Value = "2024 Annual"
$ws.Range("C2").Value = "2023 Annual"
$ws.Range("D2").Value = "2022 Annual"
